# Helper: set a cell's value as plain text, preserving the cell's original
# style (Excel auto-converts strings that look like numbers, e.g. "598.54",
# into numeric values and assigns a "Text" number-format style in the
# process; we force text semantics and then restore the original style so
# the resulting file matches the source, which stores these as plain text
# with no explicit per-cell number format).
function Set-TextValue {
    param($range, [string]$text)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") "65.615.71"
Set-TextValue $ws.Range("E2") "  +0.89%  "

Set-TextValue $ws.Range("D3") "3.186.85"
Set-TextValue $ws.Range("E3") "  +0.62%  "

Set-TextValue $ws.Range("E4") "  -0.01%  "

Set-TextValue $ws.Range("D5") "598.54"
Set-TextValue $ws.Range("E5") "  +4.65%  "

Set-TextValue $ws.Range("D6") "150.79"
Set-TextValue $ws.Range("E6") "  +0.24%  "

Set-TextValue $ws.Range("E7") "  +0.01%  "

Set-TextValue $ws.Range("D8") "3.186.22"
Set-TextValue $ws.Range("E8") "  +0.66%  "

Set-TextValue $ws.Range("D9") "0.537"
Set-TextValue $ws.Range("E9") "  +1.97%  "

Set-TextValue $ws.Range("E10") "  -1.19%  "

Set-TextValue $ws.Range("D11") "6.16"
Set-TextValue $ws.Range("E11") "  -0.79%  "

Set-TextValue $ws.Range("D12") "0.508"
Set-TextValue $ws.Range("E12") "  +0.65%  "

Set-TextValue $ws.Range("E13") "  -0.98%  "

Set-TextValue $ws.Range("D14") "38.23"
Set-TextValue $ws.Range("E14") "  +0.54%  "

Set-TextValue $ws.Range("D15") "3.715.43"
Set-TextValue $ws.Range("E15") "  +0.87%  "

Set-TextValue $ws.Range("D16") "65.823.26"
Set-TextValue $ws.Range("E16") "  +1.08%  "

Set-TextValue $ws.Range("D17") "7.29"
Set-TextValue $ws.Range("E17") "  +1.57%  "

Set-TextValue $ws.Range("D18") "3.188.33"
Set-TextValue $ws.Range("E18") "  +0.45%  "

Set-TextValue $ws.Range("E19") "  +0.32%  "

Set-TextValue $ws.Range("D20") "510.14"
Set-TextValue $ws.Range("E20") "  -0.04%  "

Set-TextValue $ws.Range("D21") "15.92"
Set-TextValue $ws.Range("E21") "  +6.79%  "

Set-TextValue $ws.Range("D22") "0.733"
Set-TextValue $ws.Range("E22") "  -0.09%  "

Set-TextValue $ws.Range("D23") "15.14"
Set-TextValue $ws.Range("E23") "  -3.44%  "

Set-TextValue $ws.Range("E24") "  +1.62%  "

Set-TextValue $ws.Range("D25") "85.28"
Set-TextValue $ws.Range("E25") "  +0.65%  "

Set-TextValue $ws.Range("E26") "  -0.13%  "

Set-TextValue $ws.Range("D27") "3.01"
Set-TextValue $ws.Range("E27") "  +3.94%  "

Set-TextValue $ws.Range("D28") "9.16"
Set-TextValue $ws.Range("E28") "  -0.05%  "

Set-TextValue $ws.Range("D29") "2.23"
Set-TextValue $ws.Range("E29") "  +1.49%  "

Set-TextValue $ws.Range("D30") "2.84"
Set-TextValue $ws.Range("E30") "  +1.31%  "

Set-TextValue $ws.Range("D31") "27.95"
Set-TextValue $ws.Range("E31") "  -0.24%  "

Set-TextValue $ws.Range("D32") "6.67"
Set-TextValue $ws.Range("E32") "  +5.83%  "

Set-TextValue $ws.Range("E33") "  +0.65%  "

Set-TextValue $ws.Range("E34") "  +0.21%  "

Set-TextValue $ws.Range("D35") "6.60"
Set-TextValue $ws.Range("E35") "  -1.03%  "

Set-TextValue $ws.Range("D36") "55.43"
Set-TextValue $ws.Range("E36") "  -0.35%  "

Set-TextValue $ws.Range("D37") "0.0915"
Set-TextValue $ws.Range("E37") "  +3.26%  "

Set-TextValue $ws.Range("D38") "483.17"
Set-TextValue $ws.Range("E38") "  +0.64%  "

Set-TextValue $ws.Range("D39") "0.0422"
Set-TextValue $ws.Range("E39") "  +0.41%  "

Set-TextValue $ws.Range("D40") "2.98"
Set-TextValue $ws.Range("E40") "  -3.77%  "

Set-TextValue $ws.Range("D41") "8.92"
Set-TextValue $ws.Range("E41") "  +3.57%  "

Set-TextValue $ws.Range("D42") "3.012.93"
Set-TextValue $ws.Range("E42") "  -3.55%  "

Set-TextValue $ws.Range("E43") "  -2.05%  "

Set-TextValue $ws.Range("D44") "0.288"
Set-TextValue $ws.Range("E44") "  -0.61%  "

Set-TextValue $ws.Range("D45") "2.43"
Set-TextValue $ws.Range("E45") "  -3.54%  "

Set-TextValue $ws.Range("E46") "  +7.45%  "

Set-TextValue $ws.Range("D47") "29.00"
Set-TextValue $ws.Range("E47") "  -0.75%  "

Set-TextValue $ws.Range("D48") "1.00"
Set-TextValue $ws.Range("E48") "  +0.00%  "

Set-TextValue $ws.Range("E49") "  +0.39%  "

Set-TextValue $ws.Range("E50") "  -0.33%  "

Set-TextValue $ws.Range("D51") "119.68"
Set-TextValue $ws.Range("E51") "  -2.35%  "
